$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.388.58'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.560.70'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.001'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '285.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.82%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3643'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '48.69'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3335'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.124'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07378'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.73'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.906'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.848'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('D16').Value = '1.561.06'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001099'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '88.65'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06721'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.305'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.98'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.49%  '
$ws.Range('D24').Value = '22.374.32'
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.386'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.559'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '149.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.37'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.012'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '122.70'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.65%  '
$ws.Range('D31').Value = '1.735.40'
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.053'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.099'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.989'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.550'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.32%  '
$ws.Range('E36').Value = '  -2.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02365'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.10%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2210'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.87%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06338'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.15%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.290'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.303'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.09'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.78%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6035'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.57'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.757'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5711'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.51'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.29%  '
$ws.Range('E49').Value = '  -4.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.207'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07214'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.64%  '
